$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, shifting existing rows 24-98 down to 25-99.
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new weekly data entry.
$ws.Cells.Item(24, 1).Value = 9
$ws.Cells.Item(24, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(24, 3).Value = "Metropolitana"
$ws.Cells.Item(24, 4).Value = 44560
$ws.Cells.Item(24, 5).Value = 13
$ws.Cells.Item(24, 6).Value = 100112022
$ws.Cells.Item(24, 7).Value = "Arveja Verde"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 43
$ws.Cells.Item(24, 11).Value = 26000
$ws.Cells.Item(24, 12).Value = 28000
$ws.Cells.Item(24, 13).Value = 27023
$ws.Cells.Item(24, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(24, 15).Value = "Carahue"
$ws.Cells.Item(24, 16).Value = 1081
$ws.Cells.Item(24, 17).Value = 25
$ws.Cells.Item(24, 18).Value = "Hortaliza"
